# Revert "Refined flow sim"
# This reverts the earlier change that added Value (in) / extra pressure
# columns (C, D) and a numeric value in B2. Restore the original two
# column "Parameters" / "Results" layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra columns C:D (Value (in), Annulus Outlet Total Pressure,
# Pintle Outlet Total Pressure, and the 0.22 value) that were added by the
# commit being reverted.
$ws.Columns("C:D").Delete()

# Clear the numeric value that was added in B2.
$ws.Range("B2").ClearContents()

# Restore the original header labels.
$ws.Range("A1").Value = "Parameters"
$ws.Range("B1").Value = "Results"

# Restore the original (non vertical-centered) bold header style and
# update the selection to match the reverted layout.
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A1:B1").VerticalAlignment = -4108

$ws.Range("B2").Select()
